$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing data values ---
# Row 13: Positionen J1 -> J3 (screw terminal connector was renumbered)
$ws.Range("D13").Value = "J3"

# Row 16: resistor value / Reichelt part number updated (220R -> 300R)
$ws.Range("C16").Value = "300R"
$ws.Range("E16").Value = "RND 0805 5 300"

# Row 19: quantity 2 -> 4, Reichelt part number now filled in
$ws.Range("A19").Value = 4
$ws.Range("E19").Value = "RND 1550805 DN"

# --- Add a new BOM row for the barrel jack (Hohlbuchse) ---
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "J"
$ws.Range("C21").Value = "Hohlbuchse"
$ws.Range("D21").Value = "J1"
$ws.Range("E21").Value = "DC-BU 072759"

# --- Grow the table (ListObject) to cover the new column/row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G21"))

# --- Add new "Geprüft" column with header + "ok" for every data row ---
$ws.Range("G1").Value = "Geprüft"
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("G$r").Value = "ok"
}

# --- Update sheet view selection to match the authored workbook ---
$ws.Range("G4").Select()
